$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The merged cells in row 7 (C7:G7, H7:K7, N7:O7 and Q7) were still using the
# default General number format; the new report values put into them are
# textual, so their style switches to Text format (numFmtId 49), same as in
# the target workbook.
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("Q7").NumberFormat = "@"

$ws.Range("C7").Value = "1 2 3 (ONE TWO THREE) 20 F.C.TABS."
$ws.Range("H7").Value = "9:0"
$ws.Range("N7").Value = "40.00"
$ws.Range("Q7").Value = "1:0"

# L7 (numFmtId 165) and P7 (numFmtId 2) keep their original number formats,
# but now hold text content, so enter them as text (leading apostrophe).
$ws.Range("L7").Value = "'1"
$ws.Range("P7").Value = "'40.0000"
